$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove trailing (non-breaking) space from the two cell values
$ws.Range("A2").Value = "CO_11_8.png"
$ws.Range("A7").Value = "Pn_105_8.png"

# Select the last cell, matching the final cursor position in the diff
$ws.Range("A16").Select()
